$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for average_doctor / average_doctor_old (BP1/BQ1)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated statistics values (Harvard case classification)
$ws.Range("E4").Value = 0.359
$ws.Range("F4").Value = 0.08699999999999999
$ws.Range("G4").Value = 0.295
$ws.Range("N4").Value = 0.405
$ws.Range("O4").Value = 0.07199999999999999
$ws.Range("P4").Value = 0.268
$ws.Range("W4").Value = 0.264
$ws.Range("X4").Value = 0.111
$ws.Range("Y4").Value = 0.333
$ws.Range("AI4").Value = 0.188
$ws.Range("AJ4").Value = 0.058
$ws.Range("AK4").Value = 0.242
$ws.Range("AU4").Value = 0.125
$ws.Range("AV4").Value = 0.018
$ws.Range("AW4").Value = 0.133
$ws.Range("BA4").Value = 1.85
$ws.Range("BB4").Value = 0.163
$ws.Range("BC4").Value = 0.403
$ws.Range("BG4").Value = 0.695
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.381
$ws.Range("BM4").Value = 0.661
$ws.Range("BN4").Value = 0.094
$ws.Range("BO4").Value = 0.306
$ws.Range("BP4").Value = 0.617
$ws.Range("BQ4").Value = 0.63
$ws.Range("E5").Value = 0.455
$ws.Range("F5").Value = 0.093
$ws.Range("G5").Value = 0.305
$ws.Range("N5").Value = 0.732
$ws.Range("O5").Value = 0.092
$ws.Range("P5").Value = 0.304
$ws.Range("W5").Value = 0.274
$ws.Range("Y5").Value = 0.353
$ws.Range("AI5").Value = 0.232
$ws.Range("AJ5").Value = 0.089
$ws.Range("AK5").Value = 0.298
$ws.Range("AU5").Value = 0.289
$ws.Range("AV5").Value = 0.095
$ws.Range("AW5").Value = 0.308
$ws.Range("BA5").Value = 1.376
$ws.Range("BB5").Value = 0.096
$ws.Range("BC5").Value = 0.31
$ws.Range("BG5").Value = 0.413
$ws.Range("BH5").Value = 0.048
$ws.Range("BI5").Value = 0.219
$ws.Range("BM5").Value = 0.579
$ws.Range("BO5").Value = 0.292
$ws.Range("BP5").Value = 0.459
$ws.Range("BQ5").Value = 0.453
$ws.Range("E6").Value = 0.401
$ws.Range("N6").Value = 0.521
$ws.Range("W6").Value = 0.269
$ws.Range("AI6").Value = 0.208
$ws.Range("AU6").Value = 0.175
$ws.Range("BA6").Value = 1.567
$ws.Range("BG6").Value = 0.518
$ws.Range("BM6").Value = 0.617
$ws.Range("BP6").Value = 0.522
$ws.Range("BQ6").Value = 0.524
$ws.Range("E7").Value = 0.432
$ws.Range("N7").Value = 0.63
$ws.Range("W7").Value = 0.272
$ws.Range("AI7").Value = 0.222
$ws.Range("AU7").Value = 0.229
$ws.Range("BA7").Value = 1.445
$ws.Range("BG7").Value = 0.449
$ws.Range("BM7").Value = 0.594
$ws.Range("BP7").Value = 0.482
$ws.Range("BQ7").Value = 0.479
$ws.Range("E8").Value = 0.479
$ws.Range("F8").Value = 0.126
$ws.Range("G8").Value = 0.354
$ws.Range("N8").Value = 0.736
$ws.Range("O8").Value = 0.073
$ws.Range("P8").Value = 0.27
$ws.Range("W8").Value = 0.265
$ws.Range("X8").Value = 0.116
$ws.Range("Y8").Value = 0.34
$ws.Range("AI8").Value = 0.216
$ws.Range("AJ8").Value = 0.094
$ws.Range("AK8").Value = 0.307
$ws.Range("AU8").Value = 0.21
$ws.Range("AV8").Value = 0.059
$ws.Range("AW8").Value = 0.243
$ws.Range("BA8").Value = 1.658
$ws.Range("BB8").Value = 0.14
$ws.Range("BC8").Value = 0.374
$ws.Range("BG8").Value = 0.546
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.327
$ws.Range("BM8").Value = 0.6870000000000001
$ws.Range("BN8").Value = 0.077
$ws.Range("BO8").Value = 0.278
$ws.Range("BP8").Value = 0.553
$ws.Range("BQ8").Value = 0.569
$ws.Range("E9").Value = 0.415
$ws.Range("F9").Value = 0.243
$ws.Range("G9").Value = 0.493
$ws.Range("N9").Value = 0.634
$ws.Range("O9").Value = 0.232
$ws.Range("P9").Value = 0.482
$ws.Range("W9").Value = 0.146
$ws.Range("X9").Value = 0.125
$ws.Range("Y9").Value = 0.353
$ws.Range("AI9").Value = 0.122
$ws.Range("AJ9").Value = 0.107
$ws.Range("AK9").Value = 0.327
$ws.Range("BA9").Value = 1.585
$ws.Range("BB9").Value = 0.238
$ws.Range("BC9").Value = 0.488
$ws.Range("BG9").Value = 0.5610000000000001
$ws.Range("BH9").Value = 0.246
$ws.Range("BI9").Value = 0.496
$ws.Range("BM9").Value = 0.634
$ws.Range("BN9").Value = 0.232
$ws.Range("BO9").Value = 0.482
$ws.Range("BP9").Value = 0.528
$ws.Range("BQ9").Value = 0.534
$ws.Range("E10").Value = 0.537
$ws.Range("N10").Value = 0.829
$ws.Range("O10").Value = 0.142
$ws.Range("P10").Value = 0.376
$ws.Range("W10").Value = 0.317
$ws.Range("X10").Value = 0.217
$ws.Range("Y10").Value = 0.465
$ws.Range("AI10").Value = 0.244
$ws.Range("AJ10").Value = 0.184
$ws.Range("AK10").Value = 0.429
$ws.Range("AU10").Value = 0.195
$ws.Range("AV10").Value = 0.157
$ws.Range("AW10").Value = 0.396
$ws.Range("BA10").Value = 1.927
$ws.Range("BB10").Value = 0.249
$ws.Range("BC10").Value = 0.499
$ws.Range("BG10").Value = 0.61
$ws.Range("BH10").Value = 0.238
$ws.Range("BI10").Value = 0.488
$ws.Range("BM10").Value = 0.854
$ws.Range("BN10").Value = 0.125
$ws.Range("BO10").Value = 0.353
$ws.Range("BP10").Value = 0.642
$ws.Range("BQ10").Value = 0.679
$ws.Range("E11").Value = 0.5610000000000001
$ws.Range("F11").Value = 0.246
$ws.Range("G11").Value = 0.496
$ws.Range("N11").Value = 0.854
$ws.Range("O11").Value = 0.125
$ws.Range("P11").Value = 0.353
$ws.Range("W11").Value = 0.317
$ws.Range("X11").Value = 0.217
$ws.Range("Y11").Value = 0.465
$ws.Range("AI11").Value = 0.244
$ws.Range("AJ11").Value = 0.184
$ws.Range("AK11").Value = 0.429
$ws.Range("AU11").Value = 0.341
$ws.Range("AV11").Value = 0.225
$ws.Range("AW11").Value = 0.474
$ws.Range("BA11").Value = 1.927
$ws.Range("BB11").Value = 0.249
$ws.Range("BC11").Value = 0.499
$ws.Range("BG11").Value = 0.61
$ws.Range("BH11").Value = 0.238
$ws.Range("BI11").Value = 0.488
$ws.Range("BM11").Value = 0.854
$ws.Range("BN11").Value = 0.125
$ws.Range("BO11").Value = 0.353
$ws.Range("BP11").Value = 0.642
$ws.Range("BQ11").Value = 0.679
$ws.Range("E12").Value = 1.435
$ws.Range("F12").Value = 0.681
$ws.Range("G12").Value = 0.825
$ws.Range("N12").Value = 1.622
$ws.Range("O12").Value = 1.587
$ws.Range("P12").Value = 1.26
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.9
$ws.Range("AJ12").Value = 0.89
$ws.Range("AU12").Value = 3
$ws.Range("AV12").Value = 1.429
$ws.Range("AW12").Value = 1.195
$ws.Range("BA12").Value = 3.714
$ws.Range("BB12").Value = 0.404
$ws.Range("BC12").Value = 0.636
$ws.Range("BG12").Value = 1.08
$ws.Range("BH12").Value = 0.074
$ws.Range("BI12").Value = 0.271
$ws.Range("BM12").Value = 1.371
$ws.Range("BN12").Value = 0.462
$ws.Range("BO12").Value = 0.68
$ws.Range("BP12").Value = 1.238
$ws.Range("BQ12").Value = 1.304
$ws.Range("E13").Value = 1.72
$ws.Range("F13").Value = 0.898
$ws.Range("G13").Value = 0.948
$ws.Range("N13").Value = 2.256
$ws.Range("O13").Value = 0.977
$ws.Range("P13").Value = 0.988
$ws.Range("W13").Value = 1.09
$ws.Range("X13").Value = 0.186
$ws.Range("Y13").Value = 0.431
$ws.Range("AI13").Value = 1.39
$ws.Range("AJ13").Value = 0.419
$ws.Range("AK13").Value = 0.647
$ws.Range("AU13").Value = 2.452
$ws.Range("AV13").Value = 0.706
$ws.Range("AW13").Value = 0.84
$ws.Range("BA13").Value = 2.606
$ws.Range("BB13").Value = 0.287
$ws.Range("BC13").Value = 0.536
$ws.Range("BG13").Value = 0.659
$ws.Range("BH13").Value = 0.08699999999999999
$ws.Range("BI13").Value = 0.294
$ws.Range("BM13").Value = 1.028
$ws.Range("BN13").Value = 0.372
$ws.Range("BO13").Value = 0.61
$ws.Range("BP13").Value = 0.869
$ws.Range("BQ13").Value = 0.799
